# Delete the last slide of the deck ("Análisis costo-beneficio",
# sldId 273 / r:id rId19 / ppt/slides/slide18.xml). All other slides,
# relationships and content stay untouched.
$p = $ppt.ActivePresentation
$lastSlide = $p.Slides.Item($p.Slides.Count)
$lastSlide.Delete()
